$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing cells (columns F and G) for various rows based on revised data
$ws.Range("F273").Value = 31674
$ws.Range("G273").Value = 1663
$ws.Range("F291").Value = 14953
$ws.Range("F306").Value = 74231
$ws.Range("G306").Value = 7487
$ws.Range("F307").Value = 77825
$ws.Range("G307").Value = 6611
$ws.Range("F310").Value = 78770
$ws.Range("G310").Value = 4086
$ws.Range("F313").Value = 73944
$ws.Range("G313").Value = 3332
$ws.Range("F314").Value = 64640
$ws.Range("G314").Value = 3211
$ws.Range("F315").Value = 57105
$ws.Range("G315").Value = 2689
$ws.Range("F316").Value = 50552
$ws.Range("G316").Value = 2276
$ws.Range("F317").Value = 63908
$ws.Range("G317").Value = 2193
$ws.Range("F320").Value = 71615
$ws.Range("G320").Value = 3260
$ws.Range("F321").Value = 94859
$ws.Range("G321").Value = 2851
$ws.Range("F322").Value = 110694
$ws.Range("G322").Value = 2398
$ws.Range("F323").Value = 217646
$ws.Range("G323").Value = 3150
$ws.Range("F324").Value = 238577
$ws.Range("G324").Value = 2758
$ws.Range("F325").Value = 757927
$ws.Range("G325").Value = 6362
$ws.Range("F326").Value = 438455
$ws.Range("G326").Value = 3929
$ws.Range("F327").Value = 235709
$ws.Range("G327").Value = 2874
$ws.Range("F328").Value = 180492
$ws.Range("G328").Value = 2655
$ws.Range("F329").Value = 83663
$ws.Range("G329").Value = 1773
$ws.Range("F330").Value = 72690
$ws.Range("G330").Value = 2096
$ws.Range("F331").Value = 152179
$ws.Range("G331").Value = 2656
$ws.Range("F332").Value = 445568
$ws.Range("G332").Value = 4448
$ws.Range("F333").Value = 270809
$ws.Range("F334").Value = 203932
$ws.Range("G334").Value = 3426
$ws.Range("F335").Value = 129594
$ws.Range("G335").Value = 2952
$ws.Range("F336").Value = 102772
$ws.Range("G336").Value = 3260
$ws.Range("F337").Value = 105538
$ws.Range("G337").Value = 2991
$ws.Range("F338").Value = 223092
$ws.Range("G338").Value = 3099
$ws.Range("F339").Value = 651627
$ws.Range("G339").Value = 5602
$ws.Range("F340").Value = 383914
$ws.Range("G340").Value = 3330
$ws.Range("F341").Value = 292944
$ws.Range("G341").Value = 3592
$ws.Range("F342").Value = 177046
$ws.Range("G342").Value = 2998
$ws.Range("F343").Value = 129150
$ws.Range("G343").Value = 2881
$ws.Range("F344").Value = 133901
$ws.Range("G344").Value = 2489
$ws.Range("F345").Value = 284379
$ws.Range("G345").Value = 3250
$ws.Range("F346").Value = 659637
$ws.Range("G346").Value = 4757
$ws.Range("F347").Value = 332801
$ws.Range("G347").Value = 2818
$ws.Range("F348").Value = 230207
$ws.Range("G348").Value = 3194
$ws.Range("F349").Value = 155380
$ws.Range("G349").Value = 2679
$ws.Range("F350").Value = 124655
$ws.Range("G350").Value = 2720
$ws.Range("F351").Value = 144817
$ws.Range("G351").Value = 2725
$ws.Range("F352").Value = 297197
$ws.Range("G352").Value = 3467
$ws.Range("F353").Value = 696633
$ws.Range("G353").Value = 5117
$ws.Range("F354").Value = 296529
$ws.Range("G354").Value = 2735
$ws.Range("F355").Value = 216795
$ws.Range("G355").Value = 3335
$ws.Range("F356").Value = 156185
$ws.Range("G356").Value = 2801
$ws.Range("F357").Value = 132560
$ws.Range("G357").Value = 2911
$ws.Range("F358").Value = 154327
$ws.Range("G358").Value = 2689

# Append two new rows of data (359 and 360)
$ws.Range("A359").Value = 44253
$ws.Range("B359").Value = 306268
$ws.Range("C359").Value = 15731
$ws.Range("D359").Value = 2848
$ws.Range("E359").Value = 7075
$ws.Range("F359").Value = 296772
$ws.Range("G359").Value = 3089

$ws.Range("A360").Value = 44254
$ws.Range("B360").Value = 308083
$ws.Range("C360").Value = 8839
$ws.Range("D360").Value = 1815
$ws.Range("E360").Value = 7189
$ws.Range("F360").Value = 615673
$ws.Range("G360").Value = 4188
